# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures (columns H:N) for a handful of leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr_ALC_46 = New-Object 'object[,]' 1,7
$arr_ALC_46[0,0] = 900
$arr_ALC_46[0,1] = 875
$arr_ALC_46[0,2] = 1000
$arr_ALC_46[0,3] = 2625
$arr_ALC_46[0,4] = 3000
$arr_ALC_46[0,5] = -2506
$arr_ALC_46[0,6] = -3238
$ws.Range("H46:N46").Value = $arr_ALC_46

$ws.Range("H59").Value = 1218.3334
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1218.3334
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = 3655.0002
$ws.Range("N59").Value = -4769.0002
$ws.Range("L59").ClearContents()

$arr_ALC_60 = New-Object 'object[,]' 1,7
$arr_ALC_60[0,0] = 900
$arr_ALC_60[0,1] = 875
$arr_ALC_60[0,2] = 1000
$arr_ALC_60[0,3] = 2625
$arr_ALC_60[0,4] = 3000
$arr_ALC_60[0,5] = -2141
$arr_ALC_60[0,6] = -3968
$ws.Range("H60:N60").Value = $arr_ALC_60

$arr_ALC_86 = New-Object 'object[,]' 1,7
$arr_ALC_86[0,0] = 2659.2
$arr_ALC_86[0,1] = 2670.2856
$arr_ALC_86[0,2] = 2633.3333
$arr_ALC_86[0,3] = 2670.2856
$arr_ALC_86[0,4] = 2633.3333
$arr_ALC_86[0,5] = -1547.2856
$arr_ALC_86[0,6] = -4879.3333
$ws.Range("H86:N86").Value = $arr_ALC_86

$arr_ALC_89 = New-Object 'object[,]' 1,7
$arr_ALC_89[0,0] = 2659.2
$arr_ALC_89[0,1] = 2670.2856
$arr_ALC_89[0,2] = 2633.3333
$arr_ALC_89[0,3] = 13351.428
$arr_ALC_89[0,4] = 13166.6665
$arr_ALC_89[0,5] = -7735.428
$arr_ALC_89[0,6] = -24398.6665
$ws.Range("H89:N89").Value = $arr_ALC_89

$arr_ALC_129 = New-Object 'object[,]' 1,7
$arr_ALC_129[0,0] = 6560.763
$arr_ALC_129[0,1] = 438
$arr_ALC_129[0,2] = 8747.464
$arr_ALC_129[0,3] = 1314
$arr_ALC_129[0,4] = 26242.392
$arr_ALC_129[0,5] = 3686
$arr_ALC_129[0,6] = -36242.392
$ws.Range("H129:N129").Value = $arr_ALC_129

$arr_ALC_137 = New-Object 'object[,]' 1,7
$arr_ALC_137[0,0] = 1234.3256
$arr_ALC_137[0,1] = 1077.4242
$arr_ALC_137[0,2] = 1752.1
$arr_ALC_137[0,3] = 3232.2726
$arr_ALC_137[0,4] = 5256.299999999999
$arr_ALC_137[0,5] = -682.2725999999998
$arr_ALC_137[0,6] = -10356.3
$ws.Range("H137:N137").Value = $arr_ALC_137

$arr_ALC_138 = New-Object 'object[,]' 1,7
$arr_ALC_138[0,0] = 2355.6924
$arr_ALC_138[0,1] = 1267.3871
$arr_ALC_138[0,2] = 3073.5107
$arr_ALC_138[0,3] = 3802.1613
$arr_ALC_138[0,4] = 9220.5321
$arr_ALC_138[0,5] = 1337.8387
$arr_ALC_138[0,6] = -19500.5321
$ws.Range("H138:N138").Value = $arr_ALC_138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 48038
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 48038
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 48038
$ws.Range("N107").Value = -55718

$arr_ARM_132 = New-Object 'object[,]' 1,7
$arr_ARM_132[0,0] = 2848.1428
$arr_ARM_132[0,1] = 2341.9412
$arr_ARM_132[0,2] = 4999.5
$arr_ARM_132[0,3] = 7025.823600000001
$arr_ARM_132[0,4] = 14998.5
$arr_ARM_132[0,5] = -4495.823600000001
$arr_ARM_132[0,6] = -20058.5
$ws.Range("H132:N132").Value = $arr_ARM_132

$ws = $wb.Worksheets.Item("BSM")
$arr_BSM_107 = New-Object 'object[,]' 1,7
$arr_BSM_107[0,0] = 1628.2354
$arr_BSM_107[0,1] = 1371.4286
$arr_BSM_107[0,2] = 2826.6667
$arr_BSM_107[0,3] = 1371.4286
$arr_BSM_107[0,4] = 2826.6667
$arr_BSM_107[0,5] = 548.5714
$arr_BSM_107[0,6] = -6666.6667
$ws.Range("H107:N107").Value = $arr_BSM_107

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 96
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 96
$ws.Range("N7").Value = -322
$ws.Range("L7").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$arr_CUL_2 = New-Object 'object[,]' 1,7
$arr_CUL_2[0,0] = 124404.375
$arr_CUL_2[0,1] = 198016.2
$arr_CUL_2[0,2] = 1718
$arr_CUL_2[0,3] = 1188097.2
$arr_CUL_2[0,4] = 10308
$arr_CUL_2[0,5] = -1187984.2
$arr_CUL_2[0,6] = -10534
$ws.Range("H2:N2").Value = $arr_CUL_2

$arr_CUL_5 = New-Object 'object[,]' 1,7
$arr_CUL_5[0,0] = 685.8511
$arr_CUL_5[0,1] = 461.74075
$arr_CUL_5[0,2] = 988.4
$arr_CUL_5[0,3] = 1385.22225
$arr_CUL_5[0,4] = 2965.2
$arr_CUL_5[0,5] = -1273.22225
$arr_CUL_5[0,6] = -3189.2
$ws.Range("H5:N5").Value = $arr_CUL_5

$arr_CUL_61 = New-Object 'object[,]' 1,7
$arr_CUL_61[0,0] = 2174
$arr_CUL_61[0,1] = 290
$arr_CUL_61[0,2] = 5000
$arr_CUL_61[0,3] = 870
$arr_CUL_61[0,4] = 15000
$arr_CUL_61[0,5] = -655
$arr_CUL_61[0,6] = -15430
$ws.Range("H61:N61").Value = $arr_CUL_61

$arr_CUL_129 = New-Object 'object[,]' 1,7
$arr_CUL_129[0,0] = 942.8889
$arr_CUL_129[0,1] = 394.75
$arr_CUL_129[0,2] = 1099.5
$arr_CUL_129[0,3] = 1184.25
$arr_CUL_129[0,4] = 3298.5
$arr_CUL_129[0,5] = 3815.75
$arr_CUL_129[0,6] = -13298.5
$ws.Range("H129:N129").Value = $arr_CUL_129

$arr_CUL_135 = New-Object 'object[,]' 1,7
$arr_CUL_135[0,0] = 685.8511
$arr_CUL_135[0,1] = 461.74075
$arr_CUL_135[0,2] = 988.4
$arr_CUL_135[0,3] = 4155.66675
$arr_CUL_135[0,4] = 8895.6
$arr_CUL_135[0,5] = -1620.66675
$arr_CUL_135[0,6] = -13965.6
$ws.Range("H135:N135").Value = $arr_CUL_135

$arr_CUL_137 = New-Object 'object[,]' 1,7
$arr_CUL_137[0,0] = 5507.643
$arr_CUL_137[0,1] = 700
$arr_CUL_137[0,2] = 6818.8184
$arr_CUL_137[0,3] = 2100
$arr_CUL_137[0,4] = 20456.4552
$arr_CUL_137[0,5] = 3000
$arr_CUL_137[0,6] = -30656.4552
$ws.Range("H137:N137").Value = $arr_CUL_137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 22828.572
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 22828.572
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 22828.572
$ws.Range("N135").Value = -32968.572

$ws = $wb.Worksheets.Item("LTW")
$arr_LTW_132 = New-Object 'object[,]' 1,7
$arr_LTW_132[0,0] = 15160454
$arr_LTW_132[0,1] = 25012708
$arr_LTW_132[0,2] = 3141
$arr_LTW_132[0,3] = 75038124
$arr_LTW_132[0,4] = 9423
$arr_LTW_132[0,5] = -75035594
$arr_LTW_132[0,6] = -14483
$ws.Range("H132:N132").Value = $arr_LTW_132

$ws.Range("H133").Value = 14919.556
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 14919.556
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 14919.556
$ws.Range("N133").Value = -19979.556

$arr_LTW_136 = New-Object 'object[,]' 1,7
$arr_LTW_136[0,0] = 10221.615
$arr_LTW_136[0,1] = 14125.111
$arr_LTW_136[0,2] = 1438.75
$arr_LTW_136[0,3] = 42375.333
$arr_LTW_136[0,4] = 4316.25
$arr_LTW_136[0,5] = -39825.333
$arr_LTW_136[0,6] = -9416.25
$ws.Range("H136:N136").Value = $arr_LTW_136

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 13768.667
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 13768.667
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 13768.667
$ws.Range("N59").Value = -15244.667

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("M61").ClearContents()

$ws.Range("H75").Value = 24565
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 24565
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 24565
$ws.Range("N75").Value = -26437

$ws.Range("H78").Value = 24565
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 24565
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 73695
$ws.Range("N78").Value = -83055

$arr_WVR_132 = New-Object 'object[,]' 1,7
$arr_WVR_132[0,0] = 2032.3928
$arr_WVR_132[0,1] = 892.0714
$arr_WVR_132[0,2] = 3172.7144
$arr_WVR_132[0,3] = 2676.2142
$arr_WVR_132[0,4] = 9518.143199999999
$arr_WVR_132[0,5] = -146.2142000000003
$arr_WVR_132[0,6] = -14578.1432
$ws.Range("H132:N132").Value = $arr_WVR_132
